$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2025-12-15 Monday" "2025-12-16 Tuesday"

Replace-Text "397÷8=49, 5" "747÷6=124, 3"
Replace-Text "331÷3=110, 1" "301÷4=75, 1"
Replace-Text "926÷7=132, 2" "659÷9=73, 2"
Replace-Text "472÷3=157, 1" "713÷4=178, 1"
Replace-Text "173÷9=19, 2" "126÷5=25, 1"

Replace-Text "532÷3=177, 1" "987÷4=246, 3"
Replace-Text "563÷7=80, 3" "529÷8=66, 1"
Replace-Text "683÷7=97, 4" "154÷9=17, 1"
Replace-Text "287÷8=35, 7" "642÷8=80, 2"
Replace-Text "167÷2=83, 1" "795÷3=265, 0"

Replace-Text "918÷8=114, 6" "553÷4=138, 1"
Replace-Text "243÷8=30, 3" "968÷9=107, 5"
Replace-Text "196÷5=39, 1" "942÷2=471, 0"
Replace-Text "853÷8=106, 5" "159÷8=19, 7"
Replace-Text "511÷2=255, 1" "515÷2=257, 1"

Replace-Text "888÷4=222, 0" "145÷2=72, 1"
Replace-Text "325÷4=81, 1" "748÷9=83, 1"
Replace-Text "475÷3=158, 1" "417÷6=69, 3"
Replace-Text "979÷6=163, 1" "160÷2=80, 0"
Replace-Text "212÷3=70, 2" "845÷6=140, 5"

Replace-Text "369÷8=46, 1" "194÷8=24, 2"
Replace-Text "414÷9=46, 0" "279÷8=34, 7"
Replace-Text "173÷7=24, 5" "278÷7=39, 5"
Replace-Text "881÷9=97, 8" "214÷2=107, 0"
Replace-Text "691÷5=138, 1" "401÷4=100, 1"
